# Apply updated cryptocurrency price/volume data to Sheet1 (cols B-E, rows 2-51).
# Uses a temporary "@" (text) number format around each assignment so that numeric-
# looking strings (e.g. "581.37", "1.00", "0.0₃0796") are stored as text, matching the
# original inline-string cells, instead of being auto-converted to numbers by Excel.
# The original cell style is restored immediately afterwards so no visible formatting
# changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '70.512.80'
Set-TextValue 'E2' '  +2.62%  '
Set-TextValue 'D3' '3.558.33'
Set-TextValue 'E3' '  +2.37%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '581.37'
Set-TextValue 'E5' '  +2.49%  '
Set-TextValue 'D6' '185.94'
Set-TextValue 'E6' '  +2.17%  '
Set-TextValue 'E7' '  +2.91%  '
Set-TextValue 'D8' '3.547.04'
Set-TextValue 'E8' '  +2.46%  '
Set-TextValue 'E9' '  -0.17%  '
Set-TextValue 'D10' '0.219'
Set-TextValue 'E10' '  +20.08%  '
Set-TextValue 'E11' '  +1.94%  '
Set-TextValue 'D12' '54.41'
Set-TextValue 'E12' '  +1.85%  '
Set-TextValue 'E13' '  +6.09%  '
Set-TextValue 'E14' '  +1.28%  '
Set-TextValue 'D15' '4.122.91'
Set-TextValue 'E15' '  +2.00%  '
Set-TextValue 'D16' '70.549.89'
Set-TextValue 'E16' '  +2.69%  '
Set-TextValue 'E17' '  +0.51%  '
Set-TextValue 'D18' '3.559.49'
Set-TextValue 'E18' '  +2.35%  '
Set-TextValue 'D19' '12.56'
Set-TextValue 'E19' '  +2.89%  '
Set-TextValue 'D20' '569.88'
Set-TextValue 'E20' '  +6.27%  '
Set-TextValue 'E21' '  +0.88%  '
Set-TextValue 'E22' '  -0.27%  '
Set-TextValue 'D23' '17.69'
Set-TextValue 'E23' '  -7.98%  '
Set-TextValue 'D24' '4.54'
Set-TextValue 'E24' '  +4.22%  '
Set-TextValue 'D25' '4.90'
Set-TextValue 'E25' '  -0.81%  '
Set-TextValue 'D26' '94.46'
Set-TextValue 'E26' '  +0.69%  '
Set-TextValue 'E27' '  +5.37%  '
Set-TextValue 'E28' '  +2.63%  '
Set-TextValue 'D29' '9.17'
Set-TextValue 'E29' '  +2.62%  '
Set-TextValue 'D30' '32.41'
Set-TextValue 'E30' '  +3.88%  '
Set-TextValue 'D31' '7.18'
Set-TextValue 'E31' '  +0.76%  '
Set-TextValue 'D32' '12.27'
Set-TextValue 'E32' '  -1.25%  '
Set-TextValue 'D33' '0.117'
Set-TextValue 'E33' '  +3.54%  '
Set-TextValue 'B34' 'Fetch.AI'
Set-TextValue 'C34' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D34' '3.43'
Set-TextValue 'E34' '  +14.47%  '
Set-TextValue 'B35' 'OKB'
Set-TextValue 'C35' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D35' '63.00'
Set-TextValue 'E35' '  -1.74%  '
Set-TextValue 'D36' '549.36'
Set-TextValue 'E36' '  -3.45%  '
Set-TextValue 'E37' '  +5.58%  '
Set-TextValue 'D38' '3.40'
Set-TextValue 'E38' '  +9.92%  '
Set-TextValue 'E39' '  +0.78%  '
Set-TextValue 'E40' '  +0.08%  '
Set-TextValue 'D41' '0.0₃0796'
Set-TextValue 'E41' '  +5.26%  '
Set-TextValue 'D42' '3.597.71'
Set-TextValue 'E42' '  +12.31%  '
Set-TextValue 'E43' '  +4.10%  '
Set-TextValue 'D44' '3.42'
Set-TextValue 'E44' '  +3.79%  '
Set-TextValue 'D45' '0.0466'
Set-TextValue 'E45' '  +7.42%  '
Set-TextValue 'D46' '3.46'
Set-TextValue 'E46' '  +0.71%  '
Set-TextValue 'E47' '  -0.70%  '
Set-TextValue 'D48' '9.30'
Set-TextValue 'E48' '  +3.58%  '
Set-TextValue 'E49' '  +3.48%  '
Set-TextValue 'D50' '1.52'
Set-TextValue 'E50' '  +16.59%  '
Set-TextValue 'D51' '1.00'
Set-TextValue 'E51' '  +0.14%  '
